# Update the "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets, which carry duplicate data tables.

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    2  = 1603
    3  = 210
    4  = 200
    5  = 3842
    6  = 6066
    7  = 341
    8  = 31
    10 = 12
    11 = 8817
    12 = 2355
    13 = 254
    14 = 5360
    15 = 10284
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
